$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q5)
$ws.Range("B7").Value = -0.2110980427227692
$ws.Range("C7").Value = 4.084665252393866
$ws.Range("D7").Value = 55.34091037395347
$ws.Range("E7").Value = 7.439147153669799
$ws.Range("F7").Value = 7.535970002723084
$ws.Range("G7").Value = 38

# Row 8 (Q6)
$ws.Range("B8").Value = -0.0736133828721533
$ws.Range("C8").Value = 4.089497862141505
$ws.Range("D8").Value = 46.25719633263552
$ws.Range("E8").Value = 6.801264318686307
$ws.Range("F8").Value = 6.894675409111837
$ws.Range("G8").Value = 37

# Row 9 (Q7)
$ws.Range("B9").Value = 0.4641329866075957
$ws.Range("C9").Value = 5.772776870582859
$ws.Range("D9").Value = 82.08122295050586
$ws.Range("E9").Value = 9.059868815303336
$ws.Range("F9").Value = 9.283023754271392
$ws.Range("G9").Value = 20

# Row 10 (Q8)
$ws.Range("B10").Value = -1.816813321953852
$ws.Range("C10").Value = 6.322456511174873
$ws.Range("D10").Value = 95.93714123768407
$ws.Range("E10").Value = 9.794750698087423
$ws.Range("F10").Value = 10.01778542427216
$ws.Range("G10").Value = 13

# Row 11 (Q9)
$ws.Range("B11").Value = -4.942193916381401
$ws.Range("C11").Value = 6.806990479524936
$ws.Range("D11").Value = 87.80047509187096
$ws.Range("E11").Value = 9.370190771370185
$ws.Range("F11").Value = 8.900505209309301
$ws.Range("G11").Value = 5
